$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.793.23'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '1.635.39'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''215.05'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = '''0.5066'
$ws.Range("E6").Value = '  -1.08%  '
$ws.Range("D7").Value = '''1.003'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '''0.2580'
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").Value = '''0.06416'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").Value = '''20.33'
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("D11").Value = '''0.07793'
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("D12").Value = '''4.248'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.640.15'
$ws.Range("E13").Value = '  +0.49%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.862.79'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '''0.5601'
$ws.Range("E15").Value = '  +1.34%  '
$ws.Range("D16").Value = '0.0₅7645'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").Value = '''63.28'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").Value = '25.810.40'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("D19").Value = '''1.002'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '''4.372'
$ws.Range("E20").Value = '  -1.27%  '
$ws.Range("D21").Value = '''192.25'
$ws.Range("E21").Value = '  -1.58%  '
$ws.Range("D22").Value = '''9.897'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").Value = '''6.146'
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '''1.775'
$ws.Range("E25").Value = '  -6.09%  '
$ws.Range("D26").Value = '''139.59'
$ws.Range("E26").Value = '  -1.95%  '
$ws.Range("D27").Value = '''0.1231'
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = '''6.812'
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("D29").Value = '''15.57'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '''1.241'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").Value = '''0.04940'
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("D32").Value = '''3.285'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("D34").Value = '''1.569'
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '''2.384'
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").Value = '''0.9017'
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("D37").Value = '''2.567'
$ws.Range("E37").Value = '  +1.08%  '
$ws.Range("D38").Value = '''0.5557'
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").Value = '1.130.53'
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("D40").Value = '''0.01569'
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").Value = '''0.9955'
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''5.447'
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''98.82'
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''0.7974'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '1.774.44'
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").Value = '0.0₈111'
$ws.Range("E46").Value = '  -5.34%  '
$ws.Range("D47").Value = '''55.60'
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("D48").Value = '''0.4251'
$ws.Range("E48").Value = '  -4.15%  '
$ws.Range("D49").Value = '''7.777'
$ws.Range("E49").Value = '  +3.00%  '
$ws.Range("E50").Value = '  -2.02%  '
$ws.Range("D51").Value = '''0.9974'
$ws.Range("E51").Value = '  -0.62%  '
